$wb = $excel.ActiveWorkbook

# --- "studies" sheet: add a new "PMID" column in H ---
$studies = $wb.Worksheets.Item("studies")
$studies.Range("H1").Value = "PMID"
$studies.Range("H2").Select()

# --- "counts" sheet: add a new "notes" column in F ---
$counts = $wb.Worksheets.Item("counts")
$counts.Range("F1").Value = "notes"
$counts.Range("F2").Select()

# Make "counts" the active / selected sheet (tabSelected + workbook activeTab)
$counts.Activate()
